$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 347 entirely (the "微笑みに微笑み。先にした人の方が美しい" post),
# which shifts all subsequent rows up by one.
$ws.Rows.Item(347).Delete()
